# Applies the "Updated cryptos list" data refresh (Sun Nov 26 23:49:44 UTC 2023).
# Only cell VALUES change (prices in column D, 1h volume % in column E, and a few
# coin rows in B/C that got re-ranked) - no formatting/style changes are intended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while keeping it stored as TEXT, even when
# the string looks like a plain number (e.g. "232.23"). Excel normally auto-converts
# such input to a numeric cell, which would lose the original text formatting, so we
# briefly force a text number-format, assign the value, then restore the default
# "Normal" style so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '37.551.11'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.066.87'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue 'D5' '232.23'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  +0.06%  '
Set-TextValue 'D8' '57.88'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('E9').Value = '  -1.53%  '
Set-TextValue 'D10' '0.0787'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '2.371.93'
$ws.Range('E13').Value = '  -0.76%  '
Set-TextValue 'D14' '21.16'
$ws.Range('E14').Value = '  -0.48%  '
Set-TextValue 'D15' '0.765'
$ws.Range('E15').Value = '  -1.85%  '
Set-TextValue 'D16' '5.34'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '2.066.83'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '37.475.23'
$ws.Range('E18').Value = '  -0.73%  '
Set-TextValue 'D19' '6.16'
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('D21').Value = '0.0₃0831'
$ws.Range('E21').Value = '  -1.87%  '
Set-TextValue 'D22' '227.72'
$ws.Range('E22').Value = '  -0.31%  '
Set-TextValue 'D23' '0.999'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +0.67%  '
Set-TextValue 'D25' '2.34'
$ws.Range('E25').Value = '  -3.19%  '
Set-TextValue 'D26' '10.04'
$ws.Range('E26').Value = '  +5.71%  '
Set-TextValue 'D27' '169.29'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('E28').Value = '  -5.20%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  -4.18%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  -3.36%  '
Set-TextValue 'D33' '0.0631'
$ws.Range('E33').Value = '  -0.91%  '
Set-TextValue 'D34' '4.68'
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('E36').Value = '  -0.11%  '
Set-TextValue 'D37' '3.32'
$ws.Range('E37').Value = '  -3.81%  '
$ws.Range('E38').Value = '  +0.10%  '
Set-TextValue 'D39' '5.29'
$ws.Range('E39').Value = '  -2.56%  '
$ws.Range('E40').Value = '  +2.98%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D41' '17.13'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '98.18'
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.489.68'
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D45' '0.0955'
$ws.Range('E45').Value = '  -2.71%  '
$ws.Range('E46').Value = '  +3.29%  '
$ws.Range('E47').Value = '  -2.50%  '
Set-TextValue 'D48' '4.01'
$ws.Range('E48').Value = '  -4.32%  '
$ws.Range('E49').Value = '  -1.34%  '
$ws.Range('E50').Value = '  -1.36%  '
$ws.Range('D51').Value = '2.256.10'
$ws.Range('E51').Value = '  -0.86%  '
